# entreculturas.xlsx data-refresh edit
# - header renames (row 1) for several columns
# - column C ('GDP') values refreshed to new model output (non-rounded floats)
# - AL column flips 0 -> 1 for a handful of rows (Colony flag)
# - the 3 cells that previously held the placeholder text '..' become numeric

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) renames ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Data refresh: column C (GDP) values + AL flag flips + '..' -> numeric ---
$ws.Range("C2").Value = 7854.952374701078
$ws.Range("C3").Value = 6923.341809163824
$ws.Range("C4").Value = 2934.187009790061
$ws.Range("C5").Value = 2870.311589353206
$ws.Range("C6").Value = 697.6889104500298
$ws.Range("AL6").Value = 1
$ws.Range("C7").Value = 1873.394108966653
$ws.Range("C8").Value = 8947.741473873051
$ws.Range("C9").Value = 1460.056109840828
$ws.Range("C10").Value = 7772.38875590225
$ws.Range("C11").Value = 5191.140356354663
$ws.Range("AL11").Value = 1
$ws.Range("C12").Value = 9502.243585046588
$ws.Range("C13").Value = 1909.084588129339
$ws.Range("C14").Value = 10594.98659239237
$ws.Range("C15").Value = 12227.21453003286
$ws.Range("C16").Value = 6128.19547247793
$ws.Range("C17").Value = 4547.50930098406
$ws.Range("C18").Value = 3972.630273980753
$ws.Range("C19").Value = 4729.735976516416
$ws.Range("C20").Value = 11155.84524560499
$ws.Range("C21").Value = 14239.03920301361
$ws.Range("C22").Value = 752.7964806390332
$ws.Range("C23").Value = 4744.762791189912
$ws.Range("C24").Value = 3928.450391496945
$ws.Range("C25").Value = 3587.883798243964
$ws.Range("C26").Value = 478.6685897045245
$ws.Range("C28").Value = 7328.615628939658
$ws.Range("C29").Value = 948.3318544592602
$ws.Range("C30").Value = 1286.515571617672
$ws.Range("C31").Value = 892.5687203369533
$ws.Range("C32").Value = 488.421401781569
$ws.Range("C33").Value = 647.8358464534491
$ws.Range("C34").Value = 1213.112645064426
$ws.Range("C35").Value = 234.2356468749991
$ws.Range("C36").Value = 951.6879611168786
$ws.Range("C37").Value = 612.3436990512633
$ws.Range("C38").Value = 1401.47747416771
$ws.Range("C39").Value = 815.8736791314819
$ws.Range("C40").Value = 2983.242707849043
$ws.Range("C41").Value = 2898.942214704482
$ws.Range("C42").Value = 665.6274194933962
$ws.Range("AL42").Value = 1
$ws.Range("C43").Value = 1904.346464968814
$ws.Range("C44").Value = 9271.398233246389
$ws.Range("C45").Value = 1503.870423231357
$ws.Range("C46").Value = 8082.02845866252
$ws.Range("C47").Value = 5555.389721901988
$ws.Range("AL47").Value = 1
$ws.Range("C48").Value = 10385.96443195552
$ws.Range("C49").Value = 1955.461557360978
$ws.Range("C50").Value = 11286.24301624575
$ws.Range("C51").Value = 12808.034586422
$ws.Range("C52").Value = 6336.709213679884
$ws.Range("C53").Value = 4633.590358399045
$ws.Range("C54").Value = 4355.934938677345
$ws.Range("C55").Value = 5082.354756663512
$ws.Range("C56").Value = 11992.01662617741
$ws.Range("C57").Value = 13825.35808833117
$ws.Range("C58").Value = 979.0516300049418
$ws.Range("C59").Value = 4094.350334420203
$ws.Range("C60").Value = 3579.960081455846
$ws.Range("C61").Value = 487.7306818514292
$ws.Range("C62").Value = 492.3430015592067
$ws.Range("C64").Value = 7454.720164783006
$ws.Range("C65").Value = 777.227218443918
$ws.Range("C66").Value = 1525.562493537689
$ws.Range("C67").Value = 1066.426265472019
$ws.Range("C68").Value = 1303.425880277445
$ws.Range("C69").Value = 863.7612548677739
$ws.Range("C70").Value = 505.2384587280311
$ws.Range("C71").Value = 2948.84548976845
$ws.Range("C72").Value = 670.2645481663891
$ws.Range("C73").Value = 1132.548400540401
$ws.Range("C74").Value = 553.2014555484933
$ws.Range("C75").Value = 235.9887666597866
$ws.Range("C76").Value = 369.2024078290272
$ws.Range("C77").Value = 982.980837581714
$ws.Range("C78").Value = 644.763840173281
$ws.Range("C79").Value = 1591.56825353313
$ws.Range("C80").Value = 864.5379000312432
$ws.Range("C81").Value = 8390.479071096475
$ws.Range("C82").Value = 3083.80337578809
$ws.Range("C83").Value = 2965.153206179127
$ws.Range("C84").Value = 691.8942672110555
$ws.Range("AL84").Value = 1
$ws.Range("C85").Value = 1939.33862702996
$ws.Range("C86").Value = 9477.887185090232
$ws.Range("C87").Value = 1577.487171555845
$ws.Range("C88").Value = 8841.561277324312
$ws.Range("C89").Value = 5660.517066940175
$ws.Range("AL89").Value = 1
$ws.Range("C90").Value = 10883.31535948899
$ws.Range("C91").Value = 2024.117324382548
$ws.Range("C92").Value = 11627.81065059172
$ws.Range("C93").Value = 13455.83781255333
$ws.Range("C94").Value = 6711.616186806423
$ws.Range("C95").Value = 4921.848409120176
$ws.Range("C96").Value = 4479.398934239905
$ws.Range("C97").Value = 5360.226632400601
$ws.Range("C98").Value = 12574.90356995006
$ws.Range("C99").Value = 14179.19231490798
$ws.Range("C100").Value = 827.7770138231788
$ws.Range("C101").Value = 4209.874800894355
$ws.Range("C102").Value = 10649.8372874572
$ws.Range("C103").Value = 2094.024217383061
$ws.Range("C104").Value = 11745.7759262897
$ws.Range("C105").Value = 14035.67913082598
$ws.Range("C106").Value = 5122.180090208862
$ws.Range("C107").Value = 4394.543881413723
$ws.Range("C108").Value = 5642.578115155247
$ws.Range("C109").Value = 12981.14011088224
$ws.Range("C110").Value = 14735.09353649063
$ws.Range("C111").Value = 6051.685746144485
$ws.Range("C112").Value = 3156.723844635973
$ws.Range("C113").Value = 2999.422762626143
$ws.Range("C114").Value = 701.4459636783288
$ws.Range("AL114").Value = 1
$ws.Range("C115").Value = 1982.009737844954
$ws.Range("C116").Value = 9690.869064532331
$ws.Range("C117").Value = 1657.651524528445
$ws.Range("C118").Value = 9541.060212621274
$ws.Range("C119").Value = 5745.422744292303
$ws.Range("AL119").Value = 1
$ws.Range("C120").Value = 873.9492833067068
$ws.Range("C121").Value = 2379.668184479739
$ws.Range("C122").Value = 4276.607903883666
$ws.Range("C123").Value = 3748.449444923865
$ws.Range("C124").Value = 482.9237812079122
$ws.Range("C125").Value = 513.7390871590731
$ws.Range("C127").Value = 7500.041066630049
$ws.Range("C128").Value = 1591.319557098113
$ws.Range("C129").Value = 1223.203431665713
$ws.Range("C130").Value = 1325.930225429421
$ws.Range("C131").Value = 909.3123437708064
$ws.Range("C132").Value = 528.6449273841434
$ws.Range("C133").Value = 1223.631935023299
$ws.Range("C134").Value = 238.8160458251716
$ws.Range("C135").Value = 1000.829216794104
$ws.Range("C136").Value = 683.460336640684
$ws.Range("C137").Value = 1745.10167474004
$ws.Range("C138").Value = 869.6014949562591
$ws.Range("C139").Value = 10784.63069920566
$ws.Range("C140").Value = 2201.396847776877
$ws.Range("C141").Value = 11993.48398487312
$ws.Range("C142").Value = 14461.17437757394
$ws.Range("C143").Value = 5295.682695961288
$ws.Range("C144").Value = 4699.493713911862
$ws.Range("C145").Value = 5919.20956823756
$ws.Range("C146").Value = 13541.20710895826
$ws.Range("C147").Value = 14721.85595470026
$ws.Range("C148").Value = 6203.843262938323
$ws.Range("C149").Value = 3212.740625904757
$ws.Range("C150").Value = 3056.152683606517
$ws.Range("C151").Value = 720.7128711178943
$ws.Range("AL151").Value = 1
$ws.Range("C152").Value = 2000.792448761861
$ws.Range("C153").Value = 9693.722968944676
$ws.Range("C154").Value = 1716.389195271215
$ws.Range("C155").Value = 10027.34623102135
$ws.Range("C156").Value = 5955.175904294275
$ws.Range("AL156").Value = 1
$ws.Range("C157").Value = 922.8902056569669
$ws.Range("C158").Value = 2497.68592515536
$ws.Range("C159").Value = 4327.37995998728
$ws.Range("C160").Value = 3796.882621798447
$ws.Range("C161").Value = 493.8183694827482
$ws.Range("C162").Value = 534.5063430177229
$ws.Range("C164").Value = 7563.992777076393
$ws.Range("C165").Value = 1620.124515672545
$ws.Range("C166").Value = 1225.558111711089
$ws.Range("C167").Value = 1360.10887014004
$ws.Range("C168").Value = 929.4690557368662
$ws.Range("C169").Value = 335.38915520098
$ws.Range("C170").Value = 1299.811672673934
$ws.Range("C171").Value = 242.8459946574492
$ws.Range("C172").Value = 1032.277326842402
$ws.Range("C173").Value = 698.3833464078615
$ws.Range("C174").Value = 1778.60982580794
$ws.Range("C175").Value = 872.1235974568563
$ws.Range("C176").Value = 4413.296891279079
$ws.Range("C177").Value = 3843.198240901342
$ws.Range("C178").Value = 10398.69400694643
$ws.Range("C179").Value = 2286.013198234259
$ws.Range("C180").Value = 11951.20944634967
$ws.Range("C181").Value = 972.7427283025324
$ws.Range("C182").Value = 1401.753174264641
$ws.Range("C183").Value = 961.3778847738438
$ws.Range("C184").Value = 14561.32616430782
$ws.Range("C185").Value = 7449.08671983612
$ws.Range("C186").Value = 1379.14068216006
$ws.Range("C187").Value = 6255.426161047989
$ws.Range("C188").Value = 7091.459432954363
$ws.Range("C189").Value = 5412.131646018807
$ws.Range("C190").Value = 3252.634165082374
$ws.Range("C191").Value = 449.4203771491282
$ws.Range("C192").Value = 3137.260298393558
$ws.Range("C193").Value = 730.3063521039821
$ws.Range("AL193").Value = 1
$ws.Range("C194").Value = 2025.814194788851
$ws.Range("C195").Value = 1060.095015975378
$ws.Range("C196").Value = 6753.607115829548
$ws.Range("C197").Value = 468.1130345750273
$ws.Range("C198").Value = 507.537974993908
$ws.Range("C199").Value = 707.8672001573369
$ws.Range("C200").Value = 3125.07948072635
$ws.Range("C201").Value = 9839.050190896
$ws.Range("C202").Value = 558.2093442539386
$ws.Range("C203").Value = 711.3043470146426
$ws.Range("C204").Value = 1775.027517189621
$ws.Range("C205").Value = 10357.504182008
$ws.Range("C206").Value = 4861.287098802361
$ws.Range("C207").Value = 5996.49696468919
$ws.Range("C208").Value = 7582.696928894958
$ws.Range("C209").Value = 612.1489724037899
$ws.Range("C210").Value = 886.4370030633224
$ws.Range("C211").Value = 14025.35756477021
$ws.Range("C212").Value = 1232.864865260161
$ws.Range("C213").Value = 831.8504623916352
$ws.Range("C214").Value = 4524.373085871202
$ws.Range("C215").Value = 3748.320622951519
$ws.Range("C216").Value = 10568.15780870825
$ws.Range("C217").Value = 2361.056581219794
$ws.Range("C218").Value = 11431.15448084494
$ws.Range("C219").Value = 1024.621364522189
$ws.Range("C220").Value = 1441.783971398429
$ws.Range("C221").Value = 956.659691840205
$ws.Range("C222").Value = 14722.36632763098
$ws.Range("C223").Value = 7580.275568826287
$ws.Range("C224").Value = 1463.71052702022
$ws.Range("C225").Value = 6522.736799041846
$ws.Range("C226").Value = 6891.120221868371
$ws.Range("C227").Value = 5330.539154475424
$ws.Range("C228").Value = 3314.741082534716
$ws.Range("C229").Value = 482.6390663355013
$ws.Range("C230").Value = 3210.869677115934
$ws.Range("C231").Value = 729.1196658666737
$ws.Range("AL231").Value = 1
$ws.Range("C232").Value = 2067.29003376698
$ws.Range("C233").Value = 1093.134170274031
$ws.Range("C234").Value = 6487.899081675427
$ws.Range("C235").Value = 469.9423670895969
$ws.Range("C236").Value = 507.5484050163182
$ws.Range("C237").Value = 729.7808175407341
$ws.Range("C238").Value = 3222.05417836739
$ws.Range("C239").Value = 10037.20149040966
$ws.Range("C240").Value = 579.0880693780265
$ws.Range("C241").Value = 731.9993357350996
$ws.Range("C242").Value = 1836.014008604312
$ws.Range("C243").Value = 10765.91029414483
$ws.Range("C244").Value = 4944.191641077407
$ws.Range("C245").Value = 6114.227214287786
$ws.Range("C246").Value = 7556.788578822353
$ws.Range("C247").Value = 630.9372503341563
$ws.Range("C248").Value = 900.3889853519216
$ws.Range("C249").Value = 0
$ws.Range("C250").Value = 1234.103352230985
$ws.Range("C251").Value = 730.9320426740553
$ws.Range("C252").Value = 3530.309422482455
$ws.Range("C253").Value = 10239.48134799327
$ws.Range("C254").Value = 2425.561644739583
$ws.Range("C255").Value = 10965.97426143915
$ws.Range("C256").Value = 1469.192636109792
$ws.Range("C257").Value = 869.0586852798759
$ws.Range("C258").Value = 14777.14884489417
$ws.Range("C259").Value = 7633.969039669125
$ws.Range("C260").Value = 6550.274372976741
$ws.Range("C261").Value = 7055.001624869326
$ws.Range("C262").Value = 5176.058803160127
$ws.Range("C263").Value = 3382.563653843273
$ws.Range("C264").Value = 3242.636921959078
$ws.Range("C265").Value = 729.8559996981501
$ws.Range("AL265").Value = 1
$ws.Range("C266").Value = 2111.193164269742
$ws.Range("C267").Value = 6411.986543373589
$ws.Range("C268").Value = 475.7454935403655
$ws.Range("C269").Value = 506.2496613373833
$ws.Range("C270").Value = 749.2194349876407
$ws.Range("C271").Value = 3212.81539531051
$ws.Range("C272").Value = 10205.79575322194
$ws.Range("C273").Value = 584.2111078769213
$ws.Range("C274").Value = 729.6614300490079
$ws.Range("C275").Value = 1895.214690888655
$ws.Range("C276").Value = 11107.22332817951
$ws.Range("C277").Value = 5089.61202008711
$ws.Range("C278").Value = 6262.368904654469
$ws.Range("C280").Value = 359.6000402964525
$ws.Range("C281").Value = 7476.621011558085
$ws.Range("C282").Value = 909.5979669529498
$ws.Range("C283").Value = 14124.14385720241
$ws.Range("C284").Value = 0
$ws.Range("C285").Value = 1224.309922120837
$ws.Range("C286").Value = 0
